# Nalco aluminium-ingot price table update.
# A brand-new day's price ("02-01-2026") is prepended to the table; every
# existing row slides down by one (row 2 -> row 3, row 3 -> row 4, ...,
# old row 149 -> new row 150), growing the sheet from A1:F149 to A1:F150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Push all data rows (2..149) down by one to make room for the new entry.
#    This also bumps the sheet dimension to A1:F150 automatically.
$ws.Rows.Item(2).Insert()

# Excel's "insert" inherits formatting from the row above (the bold header),
# so re-apply the plain data-row formatting (style used by every other row)
# by copying the (now-shifted) old row 2 -- which now lives at row 3.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 2) Fill in the new row 2 with the newest circular's data. The Date /
#    Circular Date / Circular Link columns look like dates or URLs, so we
#    briefly force text formatting while assigning them to stop Excel from
#    auto-converting the strings into real dates, then restore the normal
#    data-row format (copied from row 3 again) so the stored style matches
#    every other row in the column.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "02-01-2026"

$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 307.25

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "01-01-2026"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-01-2026.pdf"

$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3) The row-insert does not renumber the per-cell hyperlinks stored on
#    column F, so the hyperlink-to-row mapping goes stale (and the new last
#    row has no hyperlink at all). Rebuild every hyperlink from scratch using
#    each cell's own (already-correct) link text as its target.
$ws.Hyperlinks.Delete()

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $linkCell = $ws.Cells.Item($r, 6)
    $target = $linkCell.Value2
    $ws.Hyperlinks.Add($linkCell, $target)
}

# `Hyperlinks.Add` unconditionally repaints the cell with Excel's built-in
# blue/underlined "Hyperlink" style, but this sheet always displayed its
# links in the same plain centered style as every other column. Re-copy the
# original column formatting (column A's, which never touches hyperlinks)
# over column F so the visual style is unchanged by the relink.
$ws.Range("A2:A" + $lastRow).Copy()
$ws.Range("F2:F" + $lastRow).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
